# Update "想去人数" (interest count) figures in column F across all four
# sheets, matching the refreshed data pull (commit: "Update gh-pages to
# output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 195
$ws.Range("F5").Value  = 939
$ws.Range("F7").Value  = 112
$ws.Range("F8").Value  = 1354
$ws.Range("F10").Value = 48
$ws.Range("F11").Value = 83
$ws.Range("F12").Value = 109
$ws.Range("F13").Value = 39
$ws.Range("F14").Value = 1267
$ws.Range("F15").Value = 368
$ws.Range("F16").Value = 439
$ws.Range("F18").Value = 138
$ws.Range("F24").Value = 156
$ws.Range("F26").Value = 77
$ws.Range("F28").Value = 899
$ws.Range("F29").Value = 34
$ws.Range("F33").Value = 218

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value  = 690
$ws.Range("F5").Value  = 528
$ws.Range("F6").Value  = 528
$ws.Range("F11").Value = 254
$ws.Range("F16").Value = 709
$ws.Range("F19").Value = 579
$ws.Range("F21").Value = 15

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value  = 1740
$ws.Range("F5").Value  = 2032
$ws.Range("F6").Value  = 2208
$ws.Range("F7").Value  = 861
$ws.Range("F8").Value  = 852
$ws.Range("F11").Value = 1005
$ws.Range("F12").Value = 186
$ws.Range("F13").Value = 46

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 1740
$ws.Range("F3").Value  = 2032
$ws.Range("F4").Value  = 2208
$ws.Range("F9").Value  = 861
$ws.Range("F10").Value = 852
$ws.Range("F12").Value = 1005
$ws.Range("F13").Value = 195
$ws.Range("F14").Value = 186
$ws.Range("F15").Value = 46
$ws.Range("F16").Value = 690
$ws.Range("F17").Value = 939
$ws.Range("F19").Value = 112
$ws.Range("F20").Value = 1354
$ws.Range("F21").Value = 528
$ws.Range("F23").Value = 48
$ws.Range("F24").Value = 83
$ws.Range("F25").Value = 109
$ws.Range("F26").Value = 39
$ws.Range("F27").Value = 1267
$ws.Range("F28").Value = 368
$ws.Range("F29").Value = 439
$ws.Range("F30").Value = 138
$ws.Range("F36").Value = 254
$ws.Range("F37").Value = 156
$ws.Range("F38").Value = 77
$ws.Range("F40").Value = 899
$ws.Range("F42").Value = 34
$ws.Range("F44").Value = 579
$ws.Range("F45").Value = 15
$ws.Range("F50").Value = 218
